$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(875).Insert()

$ws.Cells.Item(875, 1).Value = 3
$ws.Cells.Item(875, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(875, 3).Value = "Coquimbo"
$ws.Cells.Item(875, 4).Value = 45075
$ws.Cells.Item(875, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(875, 5).Value = 5
$ws.Cells.Item(875, 6).Value = 100112045
$ws.Cells.Item(875, 7).Value = "Zapallo"
$ws.Cells.Item(875, 8).Value = "Camote"
$ws.Cells.Item(875, 9).Value = "1a (guarda)"
$ws.Cells.Item(875, 10).Value = 210
$ws.Cells.Item(875, 11).Value = 430
$ws.Cells.Item(875, 12).Value = 450
$ws.Cells.Item(875, 13).Value = 440
$ws.Cells.Item(875, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(875, 15).Value = "Provincia de Talca"
$ws.Cells.Item(875, 16).Value = 440
$ws.Cells.Item(875, 17).Value = 1
$ws.Cells.Item(875, 18).Value = "Hortaliza"
